{"js": "// Replace the 100 equation strings in the 20x5 results table (row-major\n// order) with their new values, per the target diff. Each table cell holds\n// exactly one paragraph/run, so we rewrite that paragraph's range text in\n// place to preserve the existing run formatting (font, size) and paragraph\n// alignment instead of touching the whole cell body.\nconst newValues = [\"12+14=26\", \"80+9=89\", \"3+21=24\", \"4+26=30\", \"74-34=40\", \"55-15=40\", \"53+35=88\", \"97-1=96\", \"2+48=50\", \"69+18=87\", \"90-39=51\", \"17+30=47\", \"61+10=71\", \"67-41=26\", \"72-19=53\", \"15+17=32\", \"82-76=6\", \"52+45=97\", \"40+52=92\", \"55-24=31\", \"37+49=86\", \"43+3=46\", \"94-57=37\", \"67-35=32\", \"46-28=18\", \"35+1=36\", \"63-41=22\", \"2+90=92\", \"97-65=32\", \"66-37=29\", \"35+30=65\", \"8+61=69\", \"98-29=69\", \"36-10=26\", \"75-65=10\", \"61-15=46\", \"69-17=52\", \"61-25=36\", \"50+5=55\", \"69+21=90\", \"37+50=87\", \"7+41=48\", \"99-81=18\", \"86-26=60\", \"96-7=89\", \"51+25=76\", \"97-97=0\", \"76-0=76\", \"52-43=9\", \"37+35=72\", \"43+40=83\", \"26+3=29\", \"58-1=57\", \"57-8=49\", \"82-74=8\", \"61-45=16\", \"53-45=8\", \"29+40=69\", \"66+14=80\", \"32-23=9\", \"1+92=93\", \"71-18=53\", \"56-38=18\", \"31+8=39\", \"81-8=73\", \"29+66=95\", \"4+24=28\", \"5+69=74\", \"54-41=13\", \"39+39=78\", \"77+9=86\", \"33+35=68\", \"73+14=87\", \"81-24=57\", \"42+4=46\", \"89-21=68\", \"19+31=50\", \"33+57=90\", \"60+8=68\", \"67-18=49\", \"18+78=96\", \"1+53=54\", \"8+21=29\", \"97-89=8\", \"19+75=94\", \"91-48=43\", \"41+17=58\", \"62-41=21\", \"81+12=93\", \"28-14=14\", \"71-5=66\", \"42+17=59\", \"24+46=70\", \"23-4=19\", \"37+54=91\", \"95-84=11\", \"97-28=69\", \"14-11=3\", \"51+5=56\", \"76-31=45\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst rowCount = table.values.length;\nconst colCount = table.values[0].length;\n\n// First pass: collect the first paragraph of every cell, in row-major order.\nconst cellParagraphs = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    cellParagraphs.push(paragraphs);\n  }\n}\nawait context.sync();\n\n// Second pass: overwrite each paragraph's text with the new value.\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const para = cellParagraphs[idx].items[0];\n    const range = para.getRange();\n    range.insertText(newValues[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the 100 equation strings in the 20x5 results table (row-major\n# order) with their new values, per the target diff. Setting Cell.Range.Text\n# rewrites only the text of the cell's existing run, preserving the run's\n# formatting (font, size) and the paragraph's alignment.\n$newValues = @(\n    \"12+14=26\",\n    \"80+9=89\",\n    \"3+21=24\",\n    \"4+26=30\",\n    \"74-34=40\",\n    \"55-15=40\",\n    \"53+35=88\",\n    \"97-1=96\",\n    \"2+48=50\",\n    \"69+18=87\",\n    \"90-39=51\",\n    \"17+30=47\",\n    \"61+10=71\",\n    \"67-41=26\",\n    \"72-19=53\",\n    \"15+17=32\",\n    \"82-76=6\",\n    \"52+45=97\",\n    \"40+52=92\",\n    \"55-24=31\",\n    \"37+49=86\",\n    \"43+3=46\",\n    \"94-57=37\",\n    \"67-35=32\",\n    \"46-28=18\",\n    \"35+1=36\",\n    \"63-41=22\",\n    \"2+90=92\",\n    \"97-65=32\",\n    \"66-37=29\",\n    \"35+30=65\",\n    \"8+61=69\",\n    \"98-29=69\",\n    \"36-10=26\",\n    \"75-65=10\",\n    \"61-15=46\",\n    \"69-17=52\",\n    \"61-25=36\",\n    \"50+5=55\",\n    \"69+21=90\",\n    \"37+50=87\",\n    \"7+41=48\",\n    \"99-81=18\",\n    \"86-26=60\",\n    \"96-7=89\",\n    \"51+25=76\",\n    \"97-97=0\",\n    \"76-0=76\",\n    \"52-43=9\",\n    \"37+35=72\",\n    \"43+40=83\",\n    \"26+3=29\",\n    \"58-1=57\",\n    \"57-8=49\",\n    \"82-74=8\",\n    \"61-45=16\",\n    \"53-45=8\",\n    \"29+40=69\",\n    \"66+14=80\",\n    \"32-23=9\",\n    \"1+92=93\",\n    \"71-18=53\",\n    \"56-38=18\",\n    \"31+8=39\",\n    \"81-8=73\",\n    \"29+66=95\",\n    \"4+24=28\",\n    \"5+69=74\",\n    \"54-41=13\",\n    \"39+39=78\",\n    \"77+9=86\",\n    \"33+35=68\",\n    \"73+14=87\",\n    \"81-24=57\",\n    \"42+4=46\",\n    \"89-21=68\",\n    \"19+31=50\",\n    \"33+57=90\",\n    \"60+8=68\",\n    \"67-18=49\",\n    \"18+78=96\",\n    \"1+53=54\",\n    \"8+21=29\",\n    \"97-89=8\",\n    \"19+75=94\",\n    \"91-48=43\",\n    \"41+17=58\",\n    \"62-41=21\",\n    \"81+12=93\",\n    \"28-14=14\",\n    \"71-5=66\",\n    \"42+17=59\",\n    \"24+46=70\",\n    \"23-4=19\",\n    \"37+54=91\",\n    \"95-84=11\",\n    \"97-28=69\",\n    \"14-11=3\",\n    \"51+5=56\",\n    \"76-31=45\"\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
